$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest date column (Jun_19, column D) - the remaining date
# columns to the right shift left by one, carrying their values/styles.
$ws.Range("D1").EntireColumn.Delete()

# Relabel the two newest date columns with this week's dates.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"

# Add the two new analyst rows for the single-stock group.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
